# phoneNumberExample.xlsx: start validating the phone-number column by
# formatting it as Text, and add a sample value below the header so the
# validation/format can be demonstrated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the header + new data cell as Text (numFmtId 49, "@") so phone
# numbers such as "+00000000000" are kept verbatim instead of being
# coerced into numbers.
$ws.Range("A1:A2").NumberFormat = "@"

# Add the example phone number under the "номер телефона" header.
$ws.Range("A2").Value = "+00000000000"

# Leave the selection where a user would land after typing the new row.
$ws.Range("A3").Select()
